# msz - field hint and error checks part 2
# Adds three new test-case rows worth of data (as new columns C-H in row 4)
# describing additional Vehicle Page hint/error checks, and extends the
# header row (H1) and sheet view/column-width formatting to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new column H mirrors the existing header cells (B1:G1) ---
$ws.Range("H1").Value = "dlgAutomobileInsurance"

# --- Row 4: append the new test-step cells in columns C through H ---
$ws.Range("C4").Value = "102_AutomobileInsurance_002_VehicleData_002_EnterNumericValuesBelowRange"
$ws.Range("D4").Value = "Vehicle Page check error hint list value ranges"
$ws.Range("E4").Value = "102_AutomobileInsurance_002_VehicleData_002_EnterNumericValuesAboveRange"
$ws.Range("F4").Value = "Vehicle Page check error hint list value ranges"
$ws.Range("G4").Value = "102_AutomobileInsurance_002_VehicleData_002_ManufacturingDateInTheFuture"
$ws.Range("H4").Value = "Vehicle Page check error hint manufacturing date in the future"

# --- Column widths, widened to fit the newly entered long text values ---
$ws.Columns.Item(5).ColumnWidth = 69.16666666666667
$ws.Columns.Item(7).ColumnWidth = 68.05338541666667
$ws.Columns.Item(8).ColumnWidth = 51.276041666666664

# --- View / selection state ---
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("G11").Select()
